# Make a separate column for TIMEVAL in variables1 (Close #265)
#
# The "Variables" sheet holds table "Table5" with columns:
#   pivot, order, variable-code, variable-type, en_variable-label, ...
# A new "timeval" column is inserted right after "variable-type" (i.e. it
# becomes the new column E), and the boolean TRUE that used to live in the
# "variable-type" cell for the "time" variable (row 4) now lives in this
# new "timeval" column instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")
$tbl = $ws.ListObjects.Item(1)

# Insert a blank worksheet column at E - this shifts every column at/after
# E (including the whole table body) one place to the right, preserving
# all existing cell values/formatting.
$ws.Columns("E").Insert()

# Grow the table definition so it covers the newly inserted column too.
$tbl.Resize($ws.Range("A1:Q5"))

# Re-assert every header's text so the table's column-name metadata is
# resynced against the (shifted) header cells, and give the new column
# its name.
$ws.Range("E1").Value = "timeval"
$ws.Range("F1").Value = "en_variable-label"
$ws.Range("G1").Value = "da_variable-label"
$ws.Range("H1").Value = "kl_variable-label"
$ws.Range("I1").Value = "en_domain"
$ws.Range("J1").Value = "da_domain"
$ws.Range("K1").Value = "kl_domain"
$ws.Range("L1").Value = "en_elimination"
$ws.Range("M1").Value = "da_elimination"
$ws.Range("N1").Value = "kl_elimination"
$ws.Range("O1").Value = "en_note"
$ws.Range("P1").Value = "da_note"
$ws.Range("Q1").Value = "kl_note"

# Row 4 is the "time" variable. It used to store "TIME" in the
# variable-type cell (D4); that marker now belongs in the new timeval
# column (E4) as a boolean TRUE, and D4 becomes blank again.
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = $true

# New column gets the same width the old variable-type column (D) has,
# without the auto "best fit" flag.
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth

# Reflect where the user was last working.
$ws.Range("D4").Select()
